$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.404.23"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +1.98%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.844.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.68%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'1.016"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +1.37%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'316.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +2.33%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'1.013"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +1.14%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.4740"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +1.74%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.3702"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +0.61%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.07453"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +1.28%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.8871"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +2.16%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'20.51"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +0.55%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'1.848.58"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +1.89%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.07405"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +4.52%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'5.491"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +2.58%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'93.31"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +1.66%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'6.591"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +1.65%  "
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = "'  +1.25%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'0.000008864"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +1.98%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  +1.24%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'14.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +0.78%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'27.416.38"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'5.343"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.25%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'10.71"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +1.56%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'2.076.60"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +1.99%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'1.910"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +0.57%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'152.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +0.90%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'18.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +1.66%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'2.176"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +0.56%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'5.288"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -0.29%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'118.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +2.12%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'0.08981"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +0.57%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'0.7618"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -0.55%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'1.177"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +1.74%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'4.571"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +1.63%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'2.951"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +1.60%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'1.014"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +1.28%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'1.107"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +1.97%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.05364"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +1.68%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.01968"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +0.47%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'3.013"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +2.21%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'7.323"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +0.85%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'2.395"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +1.92%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.5358"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "'0.1667"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.18%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'8.544"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +1.62%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.4962"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +0.83%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'10.51"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.62%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'1.014"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +1.28%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'105.10"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +1.38%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'1.684"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +1.05%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.06333"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +0.93%  "
$ws.Range("E51").ClearFormats()
